# Updates the crypto price-tracker sheet ("Updated cryptos list ... with GitHub Actions").
# Column D (Price) and E (Volume 1h) are refreshed with the latest scraped figures for
# every coin row; two coin pairs (rows 16/17 and 48/49) also swapped position in the
# source ranking, so their Coin name / Link / Price / Volume cells are updated together.
#
# Price/Volume are stored as literal text (e.g. "66.111.72", "  -0.63%  ") rather than
# numbers in this sheet, so plain numeric-looking prices (e.g. "585.54") are written with
# a temporary Text number format to stop Excel from auto-converting them to a float; the
# cell's original style is restored immediately afterwards so formatting is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that already carries the sheet's default (unstyled) look-and-feel, used below to
# hand a cell's style back to normal after the temporary Text format trick.
$normalStyle = $ws.Cells.Item(2, 2).Style

$updates = @(
    @{ Row = 2; Col = 4; Value = '66.111.72'; ForceText = $false }
    @{ Row = 2; Col = 5; Value = '  -0.63%  '; ForceText = $false }
    @{ Row = 3; Col = 4; Value = '3.300.78'; ForceText = $false }
    @{ Row = 3; Col = 5; Value = '  -0.29%  '; ForceText = $false }
    @{ Row = 4; Col = 5; Value = '  -0.01%  '; ForceText = $false }
    @{ Row = 5; Col = 4; Value = '585.54'; ForceText = $true }
    @{ Row = 5; Col = 5; Value = '  +2.41%  '; ForceText = $false }
    @{ Row = 6; Col = 4; Value = '181.65'; ForceText = $true }
    @{ Row = 6; Col = 5; Value = '  -0.20%  '; ForceText = $false }
    @{ Row = 7; Col = 4; Value = '0.642'; ForceText = $true }
    @{ Row = 7; Col = 5; Value = '  +7.51%  '; ForceText = $false }
    @{ Row = 8; Col = 5; Value = '  +0.00%  '; ForceText = $false }
    @{ Row = 9; Col = 5; Value = '  -2.80%  '; ForceText = $false }
    @{ Row = 10; Col = 5; Value = '  +2.41%  '; ForceText = $false }
    @{ Row = 11; Col = 4; Value = '0.406'; ForceText = $true }
    @{ Row = 11; Col = 5; Value = '  +0.90%  '; ForceText = $false }
    @{ Row = 12; Col = 4; Value = '3.871.33'; ForceText = $false }
    @{ Row = 12; Col = 5; Value = '  -0.37%  '; ForceText = $false }
    @{ Row = 13; Col = 4; Value = '0.132'; ForceText = $true }
    @{ Row = 13; Col = 5; Value = '  -4.45%  '; ForceText = $false }
    @{ Row = 14; Col = 4; Value = '66.158.56'; ForceText = $false }
    @{ Row = 14; Col = 5; Value = '  -0.69%  '; ForceText = $false }
    @{ Row = 15; Col = 4; Value = '26.33'; ForceText = $true }
    @{ Row = 15; Col = 5; Value = '  -2.85%  '; ForceText = $false }
    @{ Row = 16; Col = 2; Value = 'ShibaInu'; ForceText = $false }
    @{ Row = 16; Col = 3; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; ForceText = $false }
    @{ Row = 16; Col = 4; Value = '0.0000163'; ForceText = $true }
    @{ Row = 16; Col = 5; Value = '  -2.27%  '; ForceText = $false }
    @{ Row = 17; Col = 2; Value = 'WrappedEther'; ForceText = $false }
    @{ Row = 17; Col = 3; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; ForceText = $false }
    @{ Row = 17; Col = 4; Value = '3.278.23'; ForceText = $false }
    @{ Row = 17; Col = 5; Value = '  -0.94%  '; ForceText = $false }
    @{ Row = 18; Col = 4; Value = '428.90'; ForceText = $true }
    @{ Row = 18; Col = 5; Value = '  -0.65%  '; ForceText = $false }
    @{ Row = 19; Col = 4; Value = '13.28'; ForceText = $true }
    @{ Row = 19; Col = 5; Value = '  -2.59%  '; ForceText = $false }
    @{ Row = 20; Col = 4; Value = '5.52'; ForceText = $true }
    @{ Row = 20; Col = 5; Value = '  -2.81%  '; ForceText = $false }
    @{ Row = 21; Col = 4; Value = '7.41'; ForceText = $true }
    @{ Row = 21; Col = 5; Value = '  -2.57%  '; ForceText = $false }
    @{ Row = 22; Col = 4; Value = '71.94'; ForceText = $true }
    @{ Row = 22; Col = 5; Value = '  -2.11%  '; ForceText = $false }
    @{ Row = 23; Col = 4; Value = '0.999'; ForceText = $true }
    @{ Row = 23; Col = 5; Value = '  -0.12%  '; ForceText = $false }
    @{ Row = 24; Col = 4; Value = '5.71'; ForceText = $true }
    @{ Row = 24; Col = 5; Value = '  +0.81%  '; ForceText = $false }
    @{ Row = 25; Col = 4; Value = '3.421.22'; ForceText = $false }
    @{ Row = 25; Col = 5; Value = '  -0.90%  '; ForceText = $false }
    @{ Row = 26; Col = 4; Value = '0.514'; ForceText = $true }
    @{ Row = 26; Col = 5; Value = '  -1.00%  '; ForceText = $false }
    @{ Row = 27; Col = 5; Value = '  +2.22%  '; ForceText = $false }
    @{ Row = 28; Col = 4; Value = '0.0000113'; ForceText = $true }
    @{ Row = 28; Col = 5; Value = '  -3.80%  '; ForceText = $false }
    @{ Row = 29; Col = 4; Value = '8.99'; ForceText = $true }
    @{ Row = 29; Col = 5; Value = '  -0.54%  '; ForceText = $false }
    @{ Row = 30; Col = 4; Value = '0.998'; ForceText = $true }
    @{ Row = 30; Col = 5; Value = '  -0.19%  '; ForceText = $false }
    @{ Row = 31; Col = 4; Value = '1.96'; ForceText = $true }
    @{ Row = 31; Col = 5; Value = '  -0.02%  '; ForceText = $false }
    @{ Row = 32; Col = 4; Value = '22.40'; ForceText = $true }
    @{ Row = 32; Col = 5; Value = '  -1.50%  '; ForceText = $false }
    @{ Row = 34; Col = 4; Value = '5.21'; ForceText = $true }
    @{ Row = 34; Col = 5; Value = '  -1.71%  '; ForceText = $false }
    @{ Row = 35; Col = 4; Value = '6.62'; ForceText = $true }
    @{ Row = 35; Col = 5; Value = '  -2.22%  '; ForceText = $false }
    @{ Row = 36; Col = 4; Value = '1.20'; ForceText = $true }
    @{ Row = 36; Col = 5; Value = '  -2.71%  '; ForceText = $false }
    @{ Row = 37; Col = 4; Value = '159.11'; ForceText = $true }
    @{ Row = 37; Col = 5; Value = '  -0.29%  '; ForceText = $false }
    @{ Row = 38; Col = 4; Value = '1.45'; ForceText = $true }
    @{ Row = 38; Col = 5; Value = '  -2.80%  '; ForceText = $false }
    @{ Row = 39; Col = 4; Value = '1.82'; ForceText = $true }
    @{ Row = 39; Col = 5; Value = '  -0.97%  '; ForceText = $false }
    @{ Row = 40; Col = 4; Value = '26.74'; ForceText = $true }
    @{ Row = 40; Col = 5; Value = '  -1.44%  '; ForceText = $false }
    @{ Row = 41; Col = 4; Value = '2.854.06'; ForceText = $false }
    @{ Row = 41; Col = 5; Value = '  +1.56%  '; ForceText = $false }
    @{ Row = 42; Col = 4; Value = '0.771'; ForceText = $true }
    @{ Row = 42; Col = 5; Value = '  -2.20%  '; ForceText = $false }
    @{ Row = 43; Col = 4; Value = '4.34'; ForceText = $true }
    @{ Row = 43; Col = 5; Value = '  -1.97%  '; ForceText = $false }
    @{ Row = 44; Col = 4; Value = '40.28'; ForceText = $true }
    @{ Row = 44; Col = 5; Value = '  +0.42%  '; ForceText = $false }
    @{ Row = 45; Col = 4; Value = '0.0664'; ForceText = $true }
    @{ Row = 45; Col = 5; Value = '  -1.63%  '; ForceText = $false }
    @{ Row = 46; Col = 4; Value = '5.98'; ForceText = $true }
    @{ Row = 46; Col = 5; Value = '  -3.27%  '; ForceText = $false }
    @{ Row = 47; Col = 4; Value = '2.32'; ForceText = $true }
    @{ Row = 47; Col = 5; Value = '  -1.51%  '; ForceText = $false }
    @{ Row = 48; Col = 2; Value = 'InjectiveProtocol'; ForceText = $false }
    @{ Row = 48; Col = 3; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = $false }
    @{ Row = 48; Col = 4; Value = '23.31'; ForceText = $true }
    @{ Row = 48; Col = 5; Value = '  -4.42%  '; ForceText = $false }
    @{ Row = 49; Col = 2; Value = 'Bittensor'; ForceText = $false }
    @{ Row = 49; Col = 3; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; ForceText = $false }
    @{ Row = 49; Col = 4; Value = '315.60'; ForceText = $true }
    @{ Row = 49; Col = 5; Value = '  -1.40%  '; ForceText = $false }
    @{ Row = 50; Col = 4; Value = '0.0269'; ForceText = $true }
    @{ Row = 50; Col = 5; Value = '  -1.01%  '; ForceText = $false }
    @{ Row = 51; Col = 5; Value = '  +3.88%  '; ForceText = $false }
)

foreach ($update in $updates) {
    $cell = $ws.Cells.Item($update.Row, $update.Col)
    if ($update.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $update.Value
        $cell.Style = $normalStyle
    } else {
        $cell.Value = $update.Value
    }
}
